# VerifyEmail page, setnewPassword page, resetpasswordPage
# Rebuilds the CreateAccountTest sheet data so it reflects the new
# verify-email / reset-password test rows, drops the old 6th row, and
# re-points the hyperlinks/styles to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateAccountTest")

# ---- 1. Clear out the old hyperlinks so we can rebuild them cleanly ----
$ws.Range("A1:E6").Hyperlinks.Delete()

# ---- 2. Rewrite the data rows (header row 1 stays as-is) ----
$ws.Range("A2").Value = "testug@asu.edu"
$ws.Range("B2").Value = "Testug"
$ws.Range("C2").Value = "Test123123123"
$ws.Range("D2").Value = "Test123123123"
$ws.Range("E2").Value = "Y"

$ws.Range("A3").Value = "testug@asu.edu"
$ws.Range("B3").Value = "testug@asu.edu"
$ws.Range("C3").Value = "test123123123"
$ws.Range("D3").Value = "Test123123123"
$ws.Range("E3").Value = "Y"

$ws.Range("A4").Value = "testug@asu.edu"
$ws.Range("B4").Value = "testug@asu.edu"
$ws.Range("C4").Value = "Testsadasda"
$ws.Range("D4").Value = "Testsadasda"
$ws.Range("E4").Value = "Y"

$ws.Range("A5").Value = "testug@test.asu.edu"
$ws.Range("B5").Value = "testug@test.asu.edu"
$ws.Range("C5").Value = "Test123123123"
$ws.Range("D5").Value = "Test123123123"
$ws.Range("E5").Value = "Y"

# ---- 3. Drop the old 6th row entirely (table is now 4 data rows) ----
$ws.Range("A6:E6").EntireRow.Delete()

# ---- 4. Re-add hyperlinks on the real email cells ----
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:testug@asu.edu") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:testug@asu.edu") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:testug@asu.edu") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:testug@asu.edu") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:testug@asu.edu") | Out-Null

# ---- 5. A5/B5 carry the same "email" look even without a live link ----
$ws.Range("A5").Style = "Hyperlink"
$ws.Range("B5").Style = "Hyperlink"

# ---- 6. B2 ("Testug") is plain black text, not an email ----
$ws.Range("A1").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("B2").Font.Color = 0

# ---- 7. Password / repassword columns (C,D rows 2-5) turn blue, no underline ----
$c2 = $ws.Range("C2")
$c2.Font.Color = 12673797
$c2.Font.Underline = -4142
$c2.Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("D5").PasteSpecial(-4122)

# ---- 8. Runmode column (E rows 2-5) matches the plain header style ----
$ws.Range("A1").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E5").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---- 9. Selection moves to B4, same as the recorded edit ----
$ws.Activate()
$ws.Range("B4").Select()
